$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# "Running all the test cases": flip the Runmode column (D) from "N" to "Y"
# for every test case row that hasn't been run yet (rows 2-44).
for ($r = 2; $r -le 44; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    if ($cell.Value2 -eq "N") {
        $cell.Value = "Y"
    }
}

# Update the selection to the top of the sheet instead of the previous
# bottom-of-list cell.
$ws.Range("D4").Select()
